# Scheduled "Updated cryptos list" refresh (GitHub Actions bot).
# Re-writes the Price (D) and Volume(1h) (E) columns for every coin row
# with freshly scraped figures, and - because the source ranking shuffled
# a bit this run - rows 48/49 swap Quant <-> NEARProtocol and row 51
# becomes Aave (replacing Tezos), each with its own Coin/Link/Price/Volume.
#
# Note: many of the scraped price strings (e.g. "1.013") look like plain
# decimals to Excel's automatic type detection, but the source data stores
# them as literal text (same as the original workbook). A leading
# apostrophe is used for those so Excel keeps them as text instead of
# silently coercing them to numbers - exactly what typing '1.013 into a
# cell in the Excel UI would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.615.78'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '1.776.94'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('D4').Value = '''1.013'
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').Value = '''337.42'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').Value = '''1.008'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').Value = '''0.3890'
$ws.Range('E7').Value = '  +2.95%  '
$ws.Range('D8').Value = '''0.3439'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = '''48.16'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').Value = '''1.157'
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('D11').Value = '''0.07526'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '''1.011'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '''22.38'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').Value = '''6.445'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '1.783.83'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '''7.168'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').Value = '''0.00001090'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '''0.06725'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').Value = '''83.76'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').Value = '''1.009'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = '''17.67'
$ws.Range('E21').Value = '  +1.96%  '
$ws.Range('D22').Value = '''6.532'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('D23').Value = '27.664.34'
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').Value = '''12.34'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Value = '''2.414'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').Value = '''1.502'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('D27').Value = '''2.496'
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').Value = '''21.13'
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('D29').Value = '''155.17'
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').Value = '''137.68'
$ws.Range('E30').Value = '  +2.89%  '
$ws.Range('D31').Value = '1.984.55'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '''6.277'
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('D33').Value = '''3.992'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').Value = '''0.08894'
$ws.Range('E34').Value = '  +2.46%  '
$ws.Range('D35').Value = '''13.02'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').Value = '''0.02451'
$ws.Range('E36').Value = '  +4.73%  '
$ws.Range('D37').Value = '''5.480'
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').Value = '''0.6895'
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('D39').Value = '''0.06467'
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('D40').Value = '''0.2232'
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('D41').Value = '''1.588'
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('D42').Value = '''1.265'
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('D43').Value = '''8.489'
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('D44').Value = '''14.55'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').Value = '''1.008'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').Value = '''0.6354'
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').Value = '''3.863'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''133.60'
$ws.Range('E48').Value = '  +3.69%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''2.144'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').Value = '''0.07472'
$ws.Range('E50').Value = '  +4.89%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''79.83'
$ws.Range('E51').Value = '  +0.55%  '
